# Apply the "06-04-25 daily report" update to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header titles: update report dates from 27.03.2025 to 06.04.2025 ---
$ws.Range("A1").Value = "Mangrove Communication  06.04.2025"
$ws.Range("A11").Value = "DAILY STOCK                         (06/04/2025) "

# --- Top summary block (rows 3-6), column C values replaced; E4/F4 cleared ---
$ws.Range("C3").Value = 15834
$ws.Range("C4").Value = 24562
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("C5").Value = 25078
$ws.Range("C6").Value = 25177
# Row 7 totals (C7,E7,F7) are formulas (=SUM(...)) and recalc automatically.

# --- Row 14 (I top up) ---
$ws.Range("C14").Value = 287526
$ws.Range("D14").Value = 160651
$ws.Range("E14").Value = 200000
# G14/H14 are formulas and recalc automatically.

# --- Row 24 (New STD sim) ---
$ws.Range("C24").Value = 21
$ws.Range("D24").ClearContents()
# G24/H24 are formulas and recalc automatically.

# --- Row 26 (Rbsp sim) ---
$ws.Range("C26").Value = 1
$ws.Range("D26").ClearContents()
$ws.Range("E26").Value = 100
# G26/H26 are formulas and recalc automatically.

# --- Row 28 (Ryze sim (470)) ---
$ws.Range("E28").Value = 10
# G28/H28 are formulas and recalc automatically.

# --- Bottom summary block ---
$ws.Range("H34").Value = 45167
$ws.Range("H35").Value = 15670
$ws.Range("H38").Value = 210500
# H33 and H40 are formulas and recalc automatically.

$excel.CalculateFullRebuild()
